$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "TC_IRinvestigationIRDetails" — new "EditDetails" (E) and "Risk ID" (F)
#    columns; Risk ID holds a plain number value.
# ---------------------------------------------------------------------------
$wsDet = $wb.Worksheets.Item("TC_IRinvestigationIRDetails")

$wsDet.Range("E1").Value = "EditDetails"
$wsDet.Range("E2").Value = "edited event"
$wsDet.Range("F1").Value = "Risk ID"
$wsDet.Range("F2").Value = 63839994

[void]$wsDet.Range("A1").Copy()
[void]$wsDet.Range("E1").PasteSpecial(-4122)
[void]$wsDet.Range("F1").PasteSpecial(-4122)
[void]$wsDet.Range("A2").Copy()
[void]$wsDet.Range("E2").PasteSpecial(-4122)
[void]$wsDet.Range("F2").PasteSpecial(-4122)

$wsDet.Columns.Item(5).ColumnWidth = 13.44140625
$wsDet.Columns.Item(6).ColumnWidth = 15.21875

[void]$wsDet.Range("F5").Select()

# ---------------------------------------------------------------------------
# 2) "TC_IRinvestigationRecordInputs" — new "Evidence Description" column (C).
# ---------------------------------------------------------------------------
$wsRec = $wb.Worksheets.Item("TC_IRinvestigationRecordInputs")

$wsRec.Range("C1").Value = "Evidence Description"
$wsRec.Range("C2").Value = "test response in Record inputs"

[void]$wsRec.Range("A1").Copy()
[void]$wsRec.Range("C1").PasteSpecial(-4122)
[void]$wsRec.Range("A2").Copy()
[void]$wsRec.Range("C2").PasteSpecial(-4122)

$wsRec.Columns.Item(3).ColumnWidth = 29.5546875

[void]$wsRec.Range("C6").Select()

# ---------------------------------------------------------------------------
# 3) "TC_IRinvestigationRequestInput" — new "edited query" column (F).
# ---------------------------------------------------------------------------
$wsReq = $wb.Worksheets.Item("TC_IRinvestigationRequestInput")

$wsReq.Range("F1").Value = "edited query"
$wsReq.Range("F2").Value = "edited query for automation"

[void]$wsReq.Range("A1").Copy()
[void]$wsReq.Range("F1").PasteSpecial(-4122)
[void]$wsReq.Range("A2").Copy()
[void]$wsReq.Range("F2").PasteSpecial(-4122)

$wsReq.Columns.Item(6).ColumnWidth = 39.6640625

[void]$wsReq.Range("E17").Select()

# ---------------------------------------------------------------------------
# 4) "Test Cases" sheet — flip a few Yes -> No, tweak a description, and add
#    a brand new test-case row for the IR-details edit/delete scenario.
# ---------------------------------------------------------------------------
$wsCases = $wb.Worksheets.Item("Test Cases")

$wsCases.Range("A28").Value = "TC_IRinvestigationIRDetailsEditDelete"
$wsCases.Range("B27").Value = "Tests the investigation functionality with filling details in IR details"
$wsCases.Range("B28").Value = "Tests the edit and delete functionality in IR details"
$wsCases.Range("C28").Value = "Yes"

$wsCases.Range("C25").Value = "No"
$wsCases.Range("C26").Value = "No"
$wsCases.Range("C27").Value = "No"

[void]$wsCases.Range("A25:C25").Copy()
[void]$wsCases.Range("A28:C28").PasteSpecial(-4122)

[void]$wsCases.Range("D25").Select()

[void]$wsCases.Activate()
